$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.708.98'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '1.889.21'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''248.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").Value = '''0.9996'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '''0.4736'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '''0.2927'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("D9").Value = '''0.06531'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.49%  '
$ws.Range("D10").Value = '''21.99'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("D11").Value = '''0.07802'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '''97.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D13").Value = '1.890.87'
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").Value = '''0.7353'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").Value = '''5.247'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.54%  '
$ws.Range("D16").Value = '''285.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.21%  '
$ws.Range("D17").Value = '30.704.45'
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("D18").Value = '''13.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.31%  '
$ws.Range("D19").Value = '''0.000007537'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("D20").Value = '''0.9995'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '2.140.46'
$ws.Range("E21").Value = '  +0.84%  '
$ws.Range("D22").Value = '''5.330'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.98%  '
$ws.Range("D23").Value = '''0.9991'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = '''6.257'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.30%  '
$ws.Range("D25").Value = '''9.234'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("D26").Value = '''164.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("D27").Value = '''18.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("E29").Value = '  -1.78%  '
$ws.Range("D30").Value = '''0.09749'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.94%  '
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").Value = '''4.311'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").Value = '''4.184'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.36%  '
$ws.Range("D34").Value = '''0.04864'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.97%  '
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").Value = '''0.6974'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  +0.44%  '
$ws.Range("E38").Value = '  +2.45%  '
$ws.Range("E39").Value = '  +2.06%  '
$ws.Range("D40").Value = '''6.373'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.74%  '
$ws.Range("D41").Value = '''76.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.93%  '
$ws.Range("D42").Value = '''2.007'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.80%  '
$ws.Range("D43").Value = '''0.4260'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.51%  '
$ws.Range("D44").Value = '''0.9996'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = '''0.8356'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("D46").Value = '''101.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("D47").Value = '''9.551'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.42%  '
$ws.Range("D48").Value = '''35.75'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.82%  '
$ws.Range("D49").Value = '''7.028'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").Value = '''920.50'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("E51").Value = '  +2.11%  '
